$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 14 and 15 were still empty (default styles s="1"/s="2"/s="3").
# Copy the number formatting (date / time styles) from the previous
# filled-in row (13) so the new rows line up with style indexes 11/12,
# then fill in the new log entries.

$ws.Range("A13:C13").Copy() | Out-Null
$ws.Range("A14:C14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A13:C13").Copy() | Out-Null
$ws.Range("A15:C15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# Row 14: 24 Oct 2017, 1h20m, "Making HUD + updating HUD script (bullets + active weapon)"
$ws.Range("A14").Value = (Get-Date -Year 2017 -Month 10 -Day 24).Date
$ws.Range("B14").Value = 0.055555555555555552
$ws.Range("C14").Value = "Making HUD + updating HUD script (bullets + active weapon)"

# Row 15: 24 Oct 2017, 20 min, "Pickups bullets"
$ws.Range("A15").Value = (Get-Date -Year 2017 -Month 10 -Day 24).Date
$ws.Range("B15").Value = 0.013888888888888888
$ws.Range("C15").Value = "Pickups bullets"
